$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Set-PlainCell($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainCell "D2" "36.920.79"
Set-PlainCell "E2" "  -1.40%  "
Set-PlainCell "D3" "1.970.59"
Set-PlainCell "E3" "  -3.59%  "
Set-PlainCell "E4" "  -0.16%  "
Set-TextCell "D5" "239.95"
Set-PlainCell "E5" "  -7.11%  "
Set-TextCell "D6" "0.594"
Set-PlainCell "E6" "  -4.68%  "
Set-PlainCell "E7" "  -0.01%  "
Set-TextCell "D8" "53.20"
Set-PlainCell "E8" "  -7.88%  "
Set-TextCell "D9" "0.365"
Set-PlainCell "E9" "  -5.80%  "
Set-TextCell "D10" "0.0743"
Set-PlainCell "E10" "  -7.68%  "
Set-TextCell "D11" "0.0983"
Set-PlainCell "E11" "  -4.97%  "
Set-PlainCell "D12" "2.258.04"
Set-PlainCell "E12" "  -4.16%  "
Set-TextCell "D13" "13.73"
Set-PlainCell "E13" "  -7.76%  "
Set-TextCell "D14" "20.47"
Set-PlainCell "E14" "  -4.52%  "
Set-TextCell "D15" "0.744"
Set-PlainCell "E15" "  -9.91%  "
Set-TextCell "D16" "4.99"
Set-PlainCell "E16" "  -7.59%  "
Set-PlainCell "D17" "1.936.49"
Set-PlainCell "E17" "  -6.15%  "
Set-PlainCell "D18" "36.739.30"
Set-PlainCell "E18" "  -1.99%  "
Set-TextCell "D19" "67.54"
Set-PlainCell "E19" "  -4.02%  "
Set-PlainCell "D20" "0.0₃0801"
Set-PlainCell "E20" "  -6.88%  "
Set-TextCell "D21" "225.99"
Set-PlainCell "E21" "  -1.58%  "
Set-TextCell "D22" "4.90"
Set-PlainCell "E22" "  -6.94%  "
Set-PlainCell "E23" "  +0.09%  "
Set-TextCell "D24" "2.36"
Set-PlainCell "E24" "  -13.52%  "
Set-TextCell "D25" "2.32"
Set-PlainCell "E25" "  -1.26%  "
Set-TextCell "D26" "161.18"
Set-PlainCell "E26" "  -1.67%  "
Set-TextCell "D27" "8.52"
Set-PlainCell "E27" "  -7.68%  "
Set-PlainCell "B28" "EthereumClassic"
Set-PlainCell "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D28" "18.91"
Set-PlainCell "E28" "  -6.06%  "
Set-PlainCell "B29" "Kaspa"
Set-PlainCell "C29" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D29" "0.124"
Set-PlainCell "E29" "  -12.33%  "
Set-TextCell "D30" "1.26"
Set-PlainCell "E30" "  -7.42%  "
Set-TextCell "D31" "0.116"
Set-PlainCell "E31" "  -4.35%  "
Set-TextCell "D32" "4.35"
Set-PlainCell "E32" "  -9.35%  "
Set-TextCell "D33" "0.0605"
Set-PlainCell "E33" "  -9.68%  "
Set-TextCell "D34" "4.20"
Set-PlainCell "E34" "  -7.22%  "
Set-TextCell "D35" "2.30"
Set-PlainCell "E35" "  -7.73%  "
Set-PlainCell "B36" "BinanceUSD"
Set-PlainCell "C36" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D36" "1.00"
Set-PlainCell "E36" "  -0.11%  "
Set-PlainCell "B37" "WEMIXToken"
Set-PlainCell "C37" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D37" "1.77"
Set-PlainCell "E37" "  -2.31%  "
Set-TextCell "D38" "3.24"
Set-PlainCell "E38" "  -7.22%  "
Set-TextCell "D39" "5.06"
Set-PlainCell "E39" "  -6.52%  "
Set-TextCell "D40" "3.02"
Set-PlainCell "E40" "  -0.44%  "
Set-PlainCell "D41" "1.401.50"
Set-PlainCell "E41" "  -0.20%  "
Set-TextCell "D42" "1.12"
Set-PlainCell "E42" "  -7.86%  "
Set-TextCell "D43" "0.0884"
Set-PlainCell "E43" "  -9.07%  "
Set-TextCell "D44" "0.0201"
Set-PlainCell "E44" "  -7.67%  "
Set-PlainCell "B45" "Aave"
Set-PlainCell "C45" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D45" "86.27"
Set-PlainCell "E45" "  -6.12%  "
Set-PlainCell "B46" "InjectiveProtocol"
Set-PlainCell "C46" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D46" "15.02"
Set-PlainCell "E46" "  -8.81%  "
Set-TextCell "D47" "0.985"
Set-PlainCell "E47" "  -6.55%  "
Set-TextCell "D48" "2.85"
Set-PlainCell "E48" "  -1.33%  "
Set-TextCell "D49" "6.57"
Set-PlainCell "E49" "  -12.30%  "
Set-PlainCell "D50" "2.151.32"
Set-PlainCell "E50" "  -4.21%  "
Set-TextCell "D51" "3.52"
Set-PlainCell "E51" "  +7.60%  "
